# Update "想去人数" (F column) counts on sheets "展览" and "全部类型"
# to reflect the latest scrape snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 65
    "F3"  = 125
    "F4"  = 2033
    "F5"  = 330
    "F6"  = 577
    "F9"  = 10483
    "F12" = 276
    "F13" = 200
    "F14" = 405
    "F15" = 7374
    "F17" = 703
    "F18" = 190
    "F19" = 62
    "F20" = 3294
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $updates.Keys) {
    $ws1.Range($addr).Value = $updates[$addr]
}

$updates4 = @{
    "F2"  = 65
    "F3"  = 125
    "F4"  = 2033
    "F5"  = 330
    "F6"  = 577
    "F12" = 10483
    "F15" = 276
    "F16" = 200
    "F17" = 405
    "F18" = 7374
    "F20" = 703
    "F21" = 190
    "F22" = 62
    "F23" = 3294
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $updates4.Keys) {
    $ws4.Range($addr).Value = $updates4[$addr]
}
